$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "292.48"
Set-TextCell $ws.Range("E2") "-6.88%"
Set-TextCell $ws.Range("D3") "40.44"
Set-TextCell $ws.Range("E3") "-0.86%"
Set-TextCell $ws.Range("D4") "5.032"
Set-TextCell $ws.Range("E4") "-2.46%"
Set-TextCell $ws.Range("D5") "0.07321"
Set-TextCell $ws.Range("E5") "-3.59%"
Set-TextCell $ws.Range("D6") "1.527"
Set-TextCell $ws.Range("E6") "-9.11%"
Set-TextCell $ws.Range("D7") "0.9299"
Set-TextCell $ws.Range("E7") "-0.04%"
Set-TextCell $ws.Range("D9") "0.1175"
Set-TextCell $ws.Range("E9") "-2.08%"
Set-TextCell $ws.Range("D10") "0.1745"
Set-TextCell $ws.Range("E10") "-4.30%"
Set-TextCell $ws.Range("D11") "0.04342"
Set-TextCell $ws.Range("E11") "4.74%"
Set-TextCell $ws.Range("D12") "0.08691"
Set-TextCell $ws.Range("E12") "-4.35%"
Set-TextCell $ws.Range("D13") "0.1054"
Set-TextCell $ws.Range("E13") "0.06%"
Set-TextCell $ws.Range("D14") "0.001271"
Set-TextCell $ws.Range("E14") "-1.86%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell $ws.Range("D15") "0.005973"
Set-TextCell $ws.Range("E15") "2.40%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws.Range("D16") "3.336"
Set-TextCell $ws.Range("E16") "0.15%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell $ws.Range("D17") "4.281"
Set-TextCell $ws.Range("E17") "-1.08%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell $ws.Range("D18") "0.3289"
Set-TextCell $ws.Range("E18") "-2.04%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell $ws.Range("D19") "7.970"
Set-TextCell $ws.Range("E19") "4.43%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell $ws.Range("D20") "0.1401"
Set-TextCell $ws.Range("E20") "4.37%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextCell $ws.Range("D21") "0.2743"
Set-TextCell $ws.Range("E21") "-3.35%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell $ws.Range("D22") "0.03935"
Set-TextCell $ws.Range("E22") "-2.29%"
Set-TextCell $ws.Range("E23") "-1.27%"
Set-TextCell $ws.Range("D24") "0.003784"
Set-TextCell $ws.Range("E24") "-4.80%"
Set-TextCell $ws.Range("D26") "0.0003726"
Set-TextCell $ws.Range("D38") "0.02279"
Set-TextCell $ws.Range("E38") "-5.54%"
Set-TextCell $ws.Range("D39") "0.05044"
Set-TextCell $ws.Range("E39") "-2.42%"
Set-TextCell $ws.Range("D40") "0.006218"
Set-TextCell $ws.Range("E40") "88.29%"
Set-TextCell $ws.Range("D41") "0.007691"
Set-TextCell $ws.Range("E41") "-0.12%"
Set-TextCell $ws.Range("D42") "0.1288"
Set-TextCell $ws.Range("E42") "-0.97%"
Set-TextCell $ws.Range("D43") "0.007328"
Set-TextCell $ws.Range("E43") "-3.75%"
Set-TextCell $ws.Range("D44") "0.008268"
Set-TextCell $ws.Range("E44") "-3.65%"
Set-TextCell $ws.Range("D45") "0.2920"
Set-TextCell $ws.Range("E45") "-13.78%"
Set-TextCell $ws.Range("D46") "0.00006281"
Set-TextCell $ws.Range("E46") "-4.75%"
Set-TextCell $ws.Range("E47") "0.04%"
Set-TextCell $ws.Range("D48") "0.03157"
Set-TextCell $ws.Range("E48") "-88.25%"
Set-TextCell $ws.Range("E49") "0.04%"
Set-TextCell $ws.Range("E50") "0.04%"
